$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.1098014276795042
$ws.Range("C2").Value = 1.001597774914103
$ws.Range("D2").Value = 3.676440339101017
$ws.Range("E2").Value = 1.917404584093044
$ws.Range("F2").Value = 1.935410742146401
$ws.Range("G2").Value = 46

$ws.Range("B3").Value = -0.04263403788362221
$ws.Range("C3").Value = 0.9298606546707162
$ws.Range("D3").Value = 2.644551444389527
$ws.Range("E3").Value = 1.626207687962865
$ws.Range("F3").Value = 1.644018222113035
$ws.Range("G3").Value = 45

$ws.Range("B4").Value = 0.02349174858992152
$ws.Range("C4").Value = 0.8547905997665299
$ws.Range("D4").Value = 2.326980562185669
$ws.Range("E4").Value = 1.525444381872269
$ws.Range("F4").Value = 1.542897175107782
$ws.Range("G4").Value = 44

$ws.Range("B5").Value = -0.01048034385804278
$ws.Range("C5").Value = 0.9300756200516419
$ws.Range("D5").Value = 2.77002958153261
$ws.Range("E5").Value = 1.664340584595776
$ws.Range("F5").Value = 1.684004220309048
$ws.Range("G5").Value = 43

$ws.Range("B6").Value = 0.06150569597062733
$ws.Range("C6").Value = 0.9485881305707295
$ws.Range("D6").Value = 2.681158187928369
$ws.Range("E6").Value = 1.6374242541041
$ws.Range("F6").Value = 1.656102977579422
$ws.Range("G6").Value = 42

$ws.Range("B7").Value = 0.01952272370762674
$ws.Range("C7").Value = 0.900502877365646
$ws.Range("D7").Value = 2.519736357723951
$ws.Range("E7").Value = 1.587367744955135
$ws.Range("F7").Value = 1.606965805954676
$ws.Range("G7").Value = 41

$ws.Range("B8").Value = 0.0832317405015639
$ws.Range("C8").Value = 0.9720704585887642
$ws.Range("D8").Value = 2.678613246716973
$ws.Range("E8").Value = 1.636646952374571
$ws.Range("F8").Value = 1.655352073199598
$ws.Range("G8").Value = 40

$ws.Range("B9").Value = 0.04039162028388746
$ws.Range("C9").Value = 0.9647072254469649
$ws.Range("D9").Value = 2.676716768498583
$ws.Range("E9").Value = 1.636067470643733
$ws.Range("F9").Value = 1.656949687439857
$ws.Range("G9").Value = 39

$ws.Range("B10").Value = 0.08220434295181167
$ws.Range("C10").Value = 0.9843083280804179
$ws.Range("D10").Value = 2.693354164668678
$ws.Range("E10").Value = 1.641144163280203
$ws.Range("F10").Value = 1.661086189779386
$ws.Range("G10").Value = 38

$ws.Range("B11").Value = 0.04603842961621714
$ws.Range("C11").Value = 0.9311122998615439
$ws.Range("D11").Value = 2.722267159893553
$ws.Range("E11").Value = 1.649929440883322
$ws.Range("F11").Value = 1.672036865348211
$ws.Range("G11").Value = 37

